# Fixed scaling problems on Table.
# Adds the missing Oct 18, 2014 second time-log entry (row 90) on Sheet1,
# which ripples through the Total Time sum (E104), the Sheet2 category
# summary table, and the pie chart that reads from it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 90: new time entry -------------------------------------------------
# Date: 10/18/2014 (same day as row 89)
$ws.Range("A90").Value = 41930

# Start Time: 6:20 PM
$ws.Range("B90").Value = 0.76388888888888884

# Stop Time: 7:02 PM
$ws.Range("C90").Value = 0.79305555555555562

# Interruption: 0 mins
$ws.Range("D90").Value = 0

# Activity: Coding (E90's shared formula already lives in the cell and will
# recompute automatically once B90:D90 are populated)
$ws.Range("F90").Value = "Coding"

# --- Selection moves to A91 as the next blank row to fill in ---------------
$ws.Range("A91").Select()
